$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2 ---
$ws.Range("J2").Value = 27.3
$ws.Range("K2").Value = 36

# --- Row 3 ---
$ws.Range("E3").Value = 1.3
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1726
$ws.Range("J3").Value = 8.3699999999999992
$ws.Range("K3").Value = 36
$ws.Range("M3").Value = 0
$ws.Range("M3").NumberFormat = "0%"
$ws.Range("N3").Value = 1730

# --- Row 4 ---
$ws.Range("E4").Value = 4.3
$ws.Range("G4").Value = 0
$ws.Range("J4").Value = 14.99
$ws.Range("K4").Value = 36
$ws.Range("M4").Value = 0
$ws.Range("M4").NumberFormat = "0%"
$ws.Range("N4").Value = 1390

# --- Row 5 ---
$ws.Range("J5").Value = 60.69
$ws.Range("K5").Value = 36
$ws.Range("M5").Value = 0
$ws.Range("M5").NumberFormat = "0%"
$ws.Range("N5").Value = 1383
$ws.Range("O5").Value = 1383

# --- Row 24 ---
$ws.Range("J24").Value = 333.02
$ws.Range("K24").Value = 36

# --- Row 30 ---
$ws.Range("M30").Value = 0
$ws.Range("M30").NumberFormat = "0%"

# --- Sheet view / selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q1").Select()

Write-Host "edits applied"
